$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header values in row 1 (next to the existing title in A1)
$ws.Range("B1").Value = "von 1 bis 10"
$ws.Range("C1").Value = "von 1 bis 3"

# Update the active selection to E7 (as recorded in the saved view state)
$ws.Range("E7").Select()
